$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in homework scores that were previously empty
$ws.Range("F5").Value = 5
$ws.Range("C7").Value = 5
$ws.Range("G7").Value = 5

# Update the view so it scrolls back up and selects K5 (as in the saved file)
$ws.Range("K5").Select()
$excel.ActiveWindow.ScrollRow = 4
